# Consumption_Forecast.xlsx update
# - Shift every timestamp in column A (rows 2-97) forward by exactly one day.
# - Replace the forecasted consumption values in column B (rows 2-97) with the
#   updated figures reflecting the addition of the Necaluxan and Adrem models.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Forecasted Consumption (MW)" values for rows 2 through 97 (96 data points,
# one per 15-minute interval of the day).
$bVals = @(
    6270,6240,6210,6180,6140,6120,6120,6110,6100,6100,
    6100,6120,6150,6170,6200,6250,6320,6410,6520,6660,
    6810,6980,7160,7340,7530,7710,7880,8050,8220,8340,
    8430,8520,8590,8600,8600,8580,8530,8480,8410,8330,
    8260,8190,8140,8080,8030,8010,8000,8000,8000,8000,
    7990,7960,7930,7910,7900,7900,7910,7910,7920,7950,
    8000,8060,8120,8180,8270,8340,8410,8500,8580,8620,
    8620,8620,8600,8580,8560,8520,8480,8410,8330,8230,
    8100,7970,7840,7700,7550,7400,7250,7080,6940,6820,
    6720,6620,6400,6330,6300,6250
)

$firstRow = 2
$lastRow = 97

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $idx = $row - $firstRow

    # Column A: push the timestamp serial number forward by one full day while
    # preserving the existing time-of-day fraction exactly.
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = [double]$aCell.Value2 + 1

    # Column B: set the newly forecasted consumption value.
    $ws.Cells.Item($row, 2).Value = $bVals[$idx]
}
